# Add week 33 ("semana 33 de 2025") data in column AJ, mirroring the
# existing weekly columns (D=week1 ... AI=week32, AJ=week33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell for the new week, copying the style of the previous header (AI1):
# bold + centered, and stored as text (like the other week-number headers),
# not as a number.
$ws.Range("AJ1").NumberFormat = "@"
$ws.Range("AJ1").Font.Bold = $true
$ws.Range("AJ1").HorizontalAlignment = -4108
$ws.Range("AJ1").Value = "33"

# Data values for column AJ, keyed by row number.
# Rows not listed here keep no value in column AJ (matching the source data,
# which leaves some rows sparse/unfilled for this week).
$values = [ordered]@{
    2  = 0
    5  = 0
    6  = 38
    7  = 0
    8  = 10
    9  = 0
    10 = 0
    11 = 0
    13 = 0
    14 = 0
    15 = 0
    23 = 0
    24 = 0
    25 = 3
    26 = 0
    29 = 2
    31 = 0
    34 = 0
    35 = 12
    36 = 1
    37 = 0
    38 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    58 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 36).Value = $values[$row]
}
